# Natmi following Dr Hou advice
# Rewrites the data rows (2-7) of Sheet1 with the updated sending/target
# cluster assignments and refreshed numeric values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("ECs","Gdnf","Gfra2","ECs", 1, 0.3333333333333333, 0.1170883333333333, 0.351265, 0.09352043543917719, 0.09352043543917718, 3, 1, 0.8258696666666667, 2.477609, 0.119018388168493, 0.119018388168493, 0.09669970282055557, 0.870297325385, 0.01113065148678648, 0.01113065148678648),
    @("ECs","Gdnf","Gfra2","FAPs", 1, 0.3333333333333333, 0.1170883333333333, 0.351265, 0.09352043543917719, 0.09352043543917718, 3, 1, 2.390723333333333, 7.17217, 0.3445338280053151, 0.344533828005315, 0.2799258105611111, 2.51933229505, 0.03222095361858365, 0.03222095361858364),
    @("ECs","Gdnf","Gfra2","sCs", 1, 0.3333333333333333, 0.1170883333333333, 0.351265, 0.09352043543917719, 0.09352043543917718, 3, 1, 3.722416, 11.167248, 0.536447783826192, 0.5364477838261918, 0.4358514854133334, 3.92266336872, 0.05016883033380707, 0.05016883033380705),
    @("sCs","Gdnf","Gfra2","ECs", 3, 1, 1.134919666666667, 3.404758999999999, 0.9064795645608229, 0.9064795645608228, 3, 1, 0.8258696666666667, 2.477609, 0.119018388168493, 0.119018388168493, 0.9372957268034444, 8.435661541230999, 0.1078877366817066, 0.1078877366817065),
    @("sCs","Gdnf","Gfra2","FAPs", 3, 1, 1.134919666666667, 3.404758999999999, 0.9064795645608229, 0.9064795645608228, 3, 1, 2.390723333333333, 7.17217, 0.3445338280053151, 0.344533828005315, 2.713278928558889, 24.41951035703, 0.3123128743867314, 0.3123128743867314),
    @("sCs","Gdnf","Gfra2","sCs", 3, 1, 1.134919666666667, 3.404758999999999, 0.9064795645608229, 0.9064795645608228, 3, 1, 3.722416, 11.167248, 0.536447783826192, 0.5364477838261918, 4.224643125914667, 38.02178813323199, 0.4862789534923849, 0.4862789534923848)
)

$rowIndex = 2
foreach ($row in $data) {
    $colIndex = 1
    foreach ($value in $row) {
        $ws.Cells.Item($rowIndex, $colIndex).Value = $value
        $colIndex++
    }
    $rowIndex++
}
